$wb = $excel.ActiveWorkbook

# Login sheet
$wsLogin = $wb.Worksheets.Item("Login")
$wsLogin.Range("G2").Value = "Success - 2020/12/19 16:58:57"
$wsLogin.Range("G3").Value = "Success - 2020/12/19 16:59:00"

# School Search sheet
$wsSchool = $wb.Worksheets.Item("School Search")
$wsSchool.Range("C2").Value = "Success - 2020/12/19 16:59:03"
$wsSchool.Range("C3").Value = "Success - 2020/12/19 16:59:06"

# Product Search sheet
$wsProduct = $wb.Worksheets.Item("Product Search")
$wsProduct.Range("K2").Value = "Success - 2020/12/19 16:59:32"
$wsProduct.Range("K3").Value = "Success - 2020/12/19 16:59:54"
$wsProduct.Range("K4").Value = "Success - 2020/12/19 17:00:15"

# Shopping Cart sheet
$wsCart = $wb.Worksheets.Item("Shopping Cart")
$wsCart.Range("G2").Value = "Success - 2020/12/19 17:00:17"
$wsCart.Range("G3").Value = "Success - 2020/12/19 17:00:17"
$wsCart.Range("G4").Value = "Success - 2020/12/19 17:00:17"

# Checkout sheet
$wsCheckout = $wb.Worksheets.Item("Checkout")
$wsCheckout.Range("P2").Value = "Success - 2020/12/19 17:00:27"
$wsCheckout.Range("P3").Value = "Success - 2020/12/19 17:00:38"
$wsCheckout.Range("P4").Value = "Success - 2020/12/19 17:00:46"

# Payment sheet
$wsPayment = $wb.Worksheets.Item("Payment")
$wsPayment.Range("C2").Value = "Success - 2020/12/19 17:00:56"
